{"js": "// Updated K means value to 37\n// 1) Collapse the two numbered \"k nearest neighbor\" result lines into a\n//    single line with the new accuracy/neighbor numbers.\n{\n  const body = context.document.body;\n\n  // Delete everything from (and including) the 2nd line break through the\n  // end of the FIRST result line (keeps the leading break that already\n  // followed \"...algorithm is \"), leaving the 2nd result-line run (with its\n  // VerbatimChar style) in place so we can just retarget its text below.\n  const junk = body.search(\n    \"\\u000b1. 71.33% with 100 nearest neighbors comparison\\u000b\",\n    { matchCase: true }\n  );\n  junk.load(\"items\");\n  await context.sync();\n  if (junk.items.length > 0) {\n    junk.items[0].delete();\n    await context.sync();\n  }\n\n  const secondLine = body.search(\n    \"2. 97.66% with 50 nearest neighbors comparison\",\n    { matchCase: true }\n  );\n  secondLine.load(\"items\");\n  await context.sync();\n  if (secondLine.items.length > 0) {\n    secondLine.items[0].insertText(\n      \"96.66% with 37 nearest neighbors comparison\",\n      \"Replace\"\n    );\n    await context.sync();\n  }\n}\n\n// 2) knn(...) call: k = 100 -> k = 37\n{\n  const body = context.document.body;\n  const kArg = body.search(\"100\", { matchCase: true });\n  kArg.load(\"items\");\n  await context.sync();\n  for (const item of kArg.items) {\n    item.load(\"text\");\n  }\n  await context.sync();\n  for (const item of kArg.items) {\n    if (item.text === \"100\") {\n      item.insertText(\"37\", \"Replace\");\n      await context.sync();\n      break;\n    }\n  }\n}\n\n// 3) Confusion matrix values.\n{\n  const body = context.document.body;\n  const row0 = body.search(\"##            0 125  31\", { matchCase: true });\n  row0.load(\"items\");\n  await context.sync();\n  if (row0.items.length > 0) {\n    row0.items[0].insertText(\"##            0 154   5\", \"Replace\");\n    await context.sync();\n  }\n\n  const row1 = body.search(\"##            1  55  89\", { matchCase: true });\n  row1.load(\"items\");\n  await context.sync();\n  if (row1.items.length > 0) {\n    row1.items[0].insertText(\"##            1   6 135\", \"Replace\");\n    await context.sync();\n  }\n}\n\n// 4) Model accuracy output value.\n{\n  const body = context.document.body;\n  const acc = body.search(\"## [1] 0.7133333\", { matchCase: true });\n  acc.load(\"items\");\n  await context.sync();\n  if (acc.items.length > 0) {\n    acc.items[0].insertText(\"## [1] 0.9633333\", \"Replace\");\n    await context.sync();\n  }\n}\n\n// 5) Comment describing the neighbor count.\n{\n  const body = context.document.body;\n  const comment = body.search(\"# model accuracy with 100 neighbors\", {\n    matchCase: true,\n  });\n  comment.load(\"items\");\n  await context.sync();\n  if (comment.items.length > 0) {\n    comment.items[0].insertText(\n      \"# model accuracy with 37 neighbors\",\n      \"Replace\"\n    );\n    await context.sync();\n  }\n}\n\n// 6) KNN explanation paragraph: drop the \"deterministic algorithm\" sentence.\n{\n  const body = context.document.body;\n  const knnText = body.search(\n    \"KNN :- K-nearest neighbor works/predicts as per the surrounding datapoints (K). It is a deterministic algorithm, if you keep the value of K and run the algorithm n times, the results will be the same. KNN is lazy execution and can be applied to non-linear solutions, due to this it provides better accuracy than logistic regression\",\n    { matchCase: true }\n  );\n  knnText.load(\"items\");\n  await context.sync();\n  if (knnText.items.length > 0) {\n    knnText.items[0].insertText(\n      \"KNN :- K-nearest neighbor works/predicts as per the surrounding datapoints (K). KNN is lazy execution and can be applied to non-linear solutions, due to this it provides better accuracy than logistic regression\",\n      \"Replace\"\n    );\n    await context.sync();\n  }\n}\n", "ps1": "# Updated K means value to 37\n$d = $word.ActiveDocument\n\n# 1) Collapse the two numbered \"k nearest neighbor\" result lines into one\n#    new line, keeping the paragraph's leading line break.\n$rng = $d.Content\n$rng.Find.Execute([char]11 + \"1. 71.33% with 100 nearest neighbors comparison\" + [char]11) | Out-Null\n$rng.Delete()\n\n$rng = $d.Content\n$rng.Find.Execute(\"2. 97.66% with 50 nearest neighbors comparison\") | Out-Null\n$rng.Text = \"96.66% with 37 nearest neighbors comparison\"\n\n# 2) knn(...) call: k = 100 -> k = 37 (narrow the found range down to just\n#    the \"100\" digits so the DecValTok run keeps its style).\n$rng = $d.Content\n$rng.Find.Execute(\"k=100)\") | Out-Null\n$sub = $d.Range($rng.Start + 2, $rng.End - 1)\n$sub.Text = \"37\"\n\n# 3) Confusion matrix values.\n$rng = $d.Content\n$rng.Find.Execute(\"##            0 125  31\") | Out-Null\n$rng.Text = \"##            0 154   5\"\n\n$rng = $d.Content\n$rng.Find.Execute(\"##            1  55  89\") | Out-Null\n$rng.Text = \"##            1   6 135\"\n\n# 4) Model accuracy output value.\n$rng = $d.Content\n$rng.Find.Execute(\"## [1] 0.7133333\") | Out-Null\n$rng.Text = \"## [1] 0.9633333\"\n\n# 5) Comment describing the neighbor count.\n$rng = $d.Content\n$rng.Find.Execute(\"# model accuracy with 100 neighbors\") | Out-Null\n$rng.Text = \"# model accuracy with 37 neighbors\"\n\n# 6) KNN explanation paragraph: drop the \"deterministic algorithm\" sentence.\n$rng = $d.Content\n$rng.Find.Execute(\"KNN :- K-nearest neighbor works/predicts as per the surrounding datapoints (K). It is a deterministic algorithm, if you keep the value of K and run the algorithm n times, the results will be the same. KNN is lazy execution and can be applied to non-linear solutions, due to this it provides better accuracy than logistic regression\") | Out-Null\n$rng.Text = \"KNN :- K-nearest neighbor works/predicts as per the surrounding datapoints (K). KNN is lazy execution and can be applied to non-linear solutions, due to this it provides better accuracy than logistic regression\"\n"}
